$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D (Price) to text so numeric-looking strings
# like "216.40" or "27.233.98" are preserved exactly as text, not coerced
# into floating point numbers that would lose formatting/precision.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.233.98"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "1.635.66"
$ws.Range("E3").Value = "  -0.84%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "216.40"
$ws.Range("E5").Value = "  -0.46%  "
$ws.Range("D6").Value = "0.522"
$ws.Range("E6").Value = "  +1.20%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("E8").Value = "  -0.24%  "
$ws.Range("E9").Value = "  -0.17%  "
$ws.Range("D10").Value = "20.39"
$ws.Range("E10").Value = "  +2.23%  "
$ws.Range("D11").Value = "0.0849"
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("D12").Value = "1.638.28"
$ws.Range("E12").Value = "  -0.52%  "
$ws.Range("D13").Value = "4.15"
$ws.Range("E13").Value = "  -0.08%  "
$ws.Range("D14").Value = "0.548"
$ws.Range("E14").Value = "  +1.16%  "
$ws.Range("D15").Value = "65.32"
$ws.Range("E15").Value = "  -3.36%  "
$ws.Range("D16").Value = "27.237.15"
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("D17").Value = "0.0₃0742"
$ws.Range("E17").Value = "  +0.43%  "
$ws.Range("D18").Value = "218.30"
$ws.Range("E18").Value = "  -0.46%  "
$ws.Range("E19").Value = "  -0.21%  "
$ws.Range("D20").Value = "7.00"
$ws.Range("E20").Value = "  +1.93%  "
$ws.Range("E21").Value = "  -0.49%  "
$ws.Range("D22").Value = "2.43"
$ws.Range("E22").Value = "  -6.12%  "
$ws.Range("D23").Value = "9.09"
$ws.Range("E23").Value = "  -1.42%  "
$ws.Range("D24").Value = "147.78"
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").Value = "7.32"
$ws.Range("E26").Value = "  -3.10%  "
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("D28").Value = "15.69"
$ws.Range("E28").Value = "  -0.62%  "
$ws.Range("D29").Value = "0.0508"
$ws.Range("E29").Value = "  -0.39%  "
$ws.Range("E30").Value = "  -0.30%  "
$ws.Range("D31").Value = "3.39"
$ws.Range("E31").Value = "  -0.56%  "
$ws.Range("E32").Value = "  -1.05%  "
$ws.Range("D33").Value = "1.336.75"
$ws.Range("E33").Value = "  +5.13%  "
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("E35").Value = "  -0.39%  "
$ws.Range("D36").Value = "0.0177"
$ws.Range("E36").Value = "  -0.64%  "
$ws.Range("D37").Value = "0.548"
$ws.Range("E37").Value = "  -0.17%  "
$ws.Range("D38").Value = "0.853"
$ws.Range("E38").Value = "  +0.56%  "
$ws.Range("E39").Value = "  -0.20%  "
$ws.Range("E40").Value = "  +1.66%  "
$ws.Range("D41").Value = "0.805"
$ws.Range("E41").Value = "  -0.41%  "
$ws.Range("D42").Value = "64.53"
$ws.Range("E42").Value = "  +3.89%  "
$ws.Range("D43").Value = "1.775.55"
$ws.Range("E43").Value = "  -0.89%  "
$ws.Range("E44").Value = "  -3.42%  "
$ws.Range("D45").Value = "90.91"
$ws.Range("E45").Value = "  -0.86%  "
$ws.Range("E46").Value = "  +0.98%  "
$ws.Range("E47").Value = "  -1.96%  "
$ws.Range("D48").Value = "0.810"
$ws.Range("E48").Value = "  +21.51%  "
$ws.Range("D49").Value = "0.0514"
$ws.Range("E49").Value = "  +0.14%  "
$ws.Range("D50").Value = "0.0991"
$ws.Range("E50").Value = "  +1.55%  "
$ws.Range("D51").Value = "7.62"
$ws.Range("E51").Value = "  -0.70%  "

# Restore default (General) formatting/style on column D so cells match
# their original unstyled appearance; the values remain text strings.
$ws.Range("D2:D51").ClearFormats()
